# Apply the 'Fix heat rate modeling syntax' value updates produced by the
# upstream model recalculation. Each worksheet is addressed by its tab name
# and the affected cells are set directly to their recalculated values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76271.06239999998
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 9300.638068405267
$ws.Range("E2").Value = 2370
$ws.Range("F2").Value = 38337.2076127313

$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 162

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("R2").Value = 50.43636363636378
$ws.Range("I3").Value = 41.6
$ws.Range("M3").Value = 23.4
$ws.Range("N3").Value = 44.8531170288747
$ws.Range("P3").Value = 52
$ws.Range("Q3").Value = 52
$ws.Range("R3").Value = 31.2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 83.2
$ws.Range("N4").Value = 81.58312417100301
$ws.Range("Q4").Value = 0

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("R2").Value = 16.63636363636378
$ws.Range("I3").Value = 41.6
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 18.8531170288747
$ws.Range("P3").Value = 23.4
$ws.Range("Q3").Value = 26
$ws.Range("R3").Value = 31.2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 59.8
$ws.Range("N4").Value = 81.58312417100301
$ws.Range("Q4").Value = 0

$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("S2").Value = 10.4
$ws.Range("T2").Value = 28.31599999999998

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 188.6909090909091
$ws.Range("C2").Value = 168.9939393939394
$ws.Range("D2").Value = 155.8626262626263
$ws.Range("E2").Value = 142.7313131313131
$ws.Range("F2").Value = 129.6
$ws.Range("G2").Value = 142.47
$ws.Range("H2").Value = 170.784
$ws.Range("I2").Value = 191.376
$ws.Range("J2").Value = 214.542
$ws.Range("K2").Value = 260.874
$ws.Range("L2").Value = 322.65
$ws.Range("M2").Value = 392.148
$ws.Range("N2").Value = 469.3679999999999
$ws.Range("O2").Value = 531.1439999999999
$ws.Range("P2").Value = 585.1979999999999
$ws.Range("Q2").Value = 631.5299999999999
$ws.Range("R2").Value = 648
$ws.Range("S2").Value = 637.4949494949495
$ws.Range("T2").Value = 608.8929292929294
$ws.Range("U2").Value = 490.7111111111111
$ws.Range("V2").Value = 392.2262626262627
$ws.Range("W2").Value = 313.4383838383839
$ws.Range("X2").Value = 260.9131313131313
$ws.Range("Y2").Value = 221.5191919191919
$ws.Range("B3").Value = 182.1252525252525
$ws.Range("C3").Value = 162.4282828282828
$ws.Range("D3").Value = 149.2969696969697
$ws.Range("E3").Value = 149.2969696969697
$ws.Range("F3").Value = 149.2969696969697
$ws.Range("G3").Value = 129.6
$ws.Range("H3").Value = 129.6
$ws.Range("I3").Value = 170.784
$ws.Range("J3").Value = 232.56
$ws.Range("K3").Value = 314.928
$ws.Range("L3").Value = 314.928
$ws.Range("M3").Value = 314.928
$ws.Range("N3").Value = 333.5925858585859
$ws.Range("O3").Value = 405.664585858586
$ws.Range("P3").Value = 428.8305858585859
$ws.Range("Q3").Value = 454.570585858586
$ws.Range("R3").Value = 485.4585858585859
$ws.Range("S3").Value = 464.4484848484849
$ws.Range("T3").Value = 333.1353535353535
$ws.Range("U3").Value = 333.1353535353535
$ws.Range("V3").Value = 333.1353535353535
$ws.Range("W3").Value = 254.3474747474747
$ws.Range("X3").Value = 254.3474747474747
$ws.Range("Y3").Value = 214.9535353535353
$ws.Range("B4").Value = 168.9939393939394
$ws.Range("C4").Value = 149.2969696969697
$ws.Range("D4").Value = 149.2969696969697
$ws.Range("E4").Value = 149.2969696969697
$ws.Range("F4").Value = 149.2969696969697
$ws.Range("G4").Value = 129.6
$ws.Range("H4").Value = 129.6
$ws.Range("I4").Value = 129.6
$ws.Range("J4").Value = 139.896
$ws.Range("K4").Value = 181.08
$ws.Range("L4").Value = 181.08
$ws.Range("M4").Value = 240.282
$ws.Range("N4").Value = 321.049292929293
$ws.Range("O4").Value = 321.049292929293
$ws.Range("P4").Value = 362.2332929292929
$ws.Range("Q4").Value = 362.2332929292929
$ws.Range("R4").Value = 372.5292929292929
$ws.Range("S4").Value = 372.5292929292929
$ws.Range("T4").Value = 241.2161616161616
$ws.Range("U4").Value = 241.2161616161616
$ws.Range("V4").Value = 241.2161616161616
$ws.Range("W4").Value = 241.2161616161616
$ws.Range("X4").Value = 241.2161616161616
$ws.Range("Y4").Value = 201.8222222222222

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("T2").Value = 2.884000000000018

$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("S2").Value = 0
